$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.757.47'
$ws.Range('E2').Value = '  +4.97%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.075.56'
$ws.Range('E3').Value = '  +3.53%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '552.91'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '140.31'
$ws.Range('E6').Value = '  +8.09%  '
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.071.07'
$ws.Range('E8').Value = '  +3.60%  '
$ws.Range('E9').Value = '  +3.98%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.29'
$ws.Range('E10').Value = '  +2.60%  '
$ws.Range('E11').Value = '  +2.97%  '
$ws.Range('E12').Value = '  +5.05%  '
$ws.Range('E13').Value = '  +5.29%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.04'
$ws.Range('E14').Value = '  +6.77%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.574.61'
$ws.Range('E15').Value = '  +3.70%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.817.30'
$ws.Range('E16').Value = '  +5.06%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.072.42'
$ws.Range('E17').Value = '  +3.58%  '
$ws.Range('E18').Value = '  -0.51%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.77'
$ws.Range('E19').Value = '  +5.23%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '484.27'
$ws.Range('E20').Value = '  +6.57%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.66'
$ws.Range('E21').Value = '  +5.63%  '
$ws.Range('E22').Value = '  +3.05%  '
$ws.Range('E23').Value = '  +7.23%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '81.70'
$ws.Range('E24').Value = '  +5.42%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.61'
$ws.Range('E25').Value = '  +8.14%  '
$ws.Range('E26').Value = '  +0.19%  '
$ws.Range('E27').Value = '  +6.08%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.03'
$ws.Range('E28').Value = '  +5.05%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.00'
$ws.Range('E29').Value = '  +10.29%  '
$ws.Range('E30').Value = '  -0.06%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '26.15'
$ws.Range('E31').Value = '  +4.74%  '
$ws.Range('E32').Value = '  +2.12%  '
$ws.Range('E33').Value = '  +8.83%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.74'
$ws.Range('E34').Value = '  +8.44%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '55.96'
$ws.Range('E35').Value = '  +1.92%  '
$ws.Range('E36').Value = '  +6.06%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '470.05'
$ws.Range('E37').Value = '  +5.19%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0820'
$ws.Range('E38').Value = '  +6.37%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.158.77'
$ws.Range('E39').Value = '  +0.18%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0398'
$ws.Range('E40').Value = '  +6.27%  '
$ws.Range('E41').Value = '  +4.09%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.29'
$ws.Range('E42').Value = '  +4.55%  '
$ws.Range('E43').Value = '  +10.90%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '28.39'
$ws.Range('E44').Value = '  +13.10%  '
$ws.Range('E45').Value = '  +5.60%  '
$ws.Range('E46').Value = '  -0.11%  '
$ws.Range('E47').Value = '  +7.85%  '
$ws.Range('E48').Value = '  +2.91%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0₃0516'
$ws.Range('E49').Value = '  +3.04%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '116.66'
$ws.Range('E50').Value = '  -0.75%  '
$ws.Range('E51').Value = '  +7.93%  '
